$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.322.09"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "3.323.59"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'586.28"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").Value = "'183.55"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "'0.646"
$ws.Range("E7").Value = "  +7.63%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.65%  "
$ws.Range("D10").Value = "'6.79"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("D11").Value = "'0.403"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "3.901.59"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("D14").Value = "66.362.62"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "'26.39"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000164"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.231.20"
$ws.Range("E17").Value = "  -3.48%  "
$ws.Range("D18").Value = "'430.22"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").Value = "'13.34"
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("D20").Value = "'5.54"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").Value = "'7.44"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").Value = "'72.09"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").Value = "3.460.76"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").Value = "'0.519"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").Value = "'0.205"
$ws.Range("E27").Value = "  +7.34%  "
$ws.Range("D29").Value = "'9.00"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "'1.94"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").Value = "'22.42"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("E36").Value = "  -4.28%  "
$ws.Range("D37").Value = "'159.87"
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("D39").Value = "2.901.07"
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("D41").Value = "'26.68"
$ws.Range("E41").Value = "  -3.50%  "
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("D44").Value = "'40.17"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("E46").Value = "  -3.68%  "
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").Value = "'23.38"
$ws.Range("E48").Value = "  -5.34%  "
$ws.Range("D49").Value = "'316.09"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("E51").Value = "  +4.69%  "

# Reset style on cells where an apostrophe-prefixed literal was used,
# so Excel does not tag them with a quotePrefix style (keeps original unstyled cells).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
